$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1, J1 and copy the style used by the existing
# header cells (e.g. H1) so formatting (bold, border, alignment) matches.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new data columns I and J for rows 2-8
$values = @(
    @(3, 6),
    @(5, 7),
    @(1, 3),
    @(7, 9),
    @(6, 8),
    @(4, 6),
    @(3, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
